$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)

# "Calculate BSdistance/EnergyConsumption/SIR/Reward" -> "Calculate BS distance, Energy consumption, SIR, Reward"
$sh1 = $s.Shapes.Item(26)
$tr1 = $sh1.TextFrame.TextRange
$tr1.Characters(1, $tr1.Length).Text = "Calculate BS distance, Energy consumption, SIR, Reward"

# "Reset parameters(state, UEposition) " -> "Reset parameters(state, UE position) "
# Only the first two runs ("Reset parameters(state, " + "UEposition") are merged;
# the trailing ")" and " " runs are left untouched.
$sh2 = $s.Shapes.Item(27)
$tr2 = $sh2.TextFrame.TextRange
$tr2.Characters(1, 34).Text = "Reset parameters(state, UE position"

# "Associate BS/UE by association rule" -> "Associate BS-UE by association rule"
$sh3 = $s.Shapes.Item(28)
$tr3 = $sh3.TextFrame.TextRange
$tr3.Characters(1, $tr3.Length).Text = "Associate BS-UE by association rule"
